$d = $word.ActiveDocument

$r1 = $d.Content
$r1.Find.Execute("Die 5. Karte", $true, $false, $false, $false, $false, $true, 1, $false, "Die 5. Karte", 2) | Out-Null

$r2 = $d.Content
$r2.Find.Execute(" ist die letzte Karte ", $true, $false, $false, $false, $false, $true, 1, $false, " ist die letzte Karte ", 2) | Out-Null
